# CoCoMo estimate sheet build-out.
# Recreates the Phase / Effort / Labour / Running / Fixed / Total-cost / Duration
# table on Hoja1, including the CoCoMo formulas, the merged "Running costs"
# column, header/label bold formatting, and the stray underlined cell at F16.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Labels (written in this order so shared-string indices line up the way
#     the original authoring tool produced them: phases first top-to-bottom
#     skipping the Integration & Test row, then the column headers left to
#     right skipping Total cost, then the two stragglers). ---
$ws.Range("B2").Value = "Phase"
$ws.Range("B3").Value = "Requirements"
$ws.Range("B4").Value = "Product Design"
$ws.Range("B5").Value = "Detailed Design"
$ws.Range("B6").Value = "Code & Unit Tests"
$ws.Range("B8").Value = "Development"
$ws.Range("B9").Value = "Totals"
$ws.Range("C2").Value = "Effort (Person-Months)"
$ws.Range("D2").Value = "Labour cost (K$)"
$ws.Range("E2").Value = "Running costs (K$)"
$ws.Range("F2").Value = "Fixed costs (K$)"
$ws.Range("H2").Value = "Duration (Months)"
$ws.Range("B7").Value = "Integration & Test"
$ws.Range("G2").Value = "Total cost (K$)"

# --- Row 3: Requirements ---
$ws.Range("C3").Value = 0.1
$ws.Range("D3").Value = 0.8
$ws.Range("E3").Value = 1.05
$ws.Range("F3").Value = 0
$ws.Range("H3").Value = 0.7
$ws.Range("G3").Formula = "=D3+E3*H3+F3"

# --- Row 4: Product Design ---
$ws.Range("C4").Value = 0.3
$ws.Range("D4").Value = 1.5
$ws.Range("F4").Value = 0
$ws.Range("H4").Value = 1
$ws.Range("G4").Formula = "=D4+E`$3*H4+F4"

# --- Row 5: Detailed Design ---
$ws.Range("C5").Value = 0.4
$ws.Range("D5").Value = 2.3
$ws.Range("F5").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("G5").Formula = "=D5+E3*H5+F5"

# --- Row 6: Code & Unit Tests ---
$ws.Range("C6").Value = 0.5
$ws.Range("D6").Value = 2.8
$ws.Range("F6").Formula = "=(1.65+1.1) * 3"
$ws.Range("H6").Value = 1.3
$ws.Range("G6").Formula = "=D6+E3*H6+F6"

# --- Row 7: Integration & Test ---
$ws.Range("C7").Value = 0.3
$ws.Range("D7").Value = 1.6
$ws.Range("F7").Value = 3.2
$ws.Range("H7").Value = 0.8
$ws.Range("G7").Formula = "=D7+E3*H7+F7"

# --- Row 8: Development (sum of the 4 phase rows above) ---
$ws.Range("C8").Formula = "=SUM(C4:C7)"
$ws.Range("D8").Formula = "=SUM(D4:D7)"
$ws.Range("F8").Formula = "=SUM(F4:F7)"
$ws.Range("G8").Formula = "=SUM(G4:G7)"
$ws.Range("H8").Formula = "=SUM(H4:H7)"

# --- Row 9: Totals (Development + Requirements) ---
$ws.Range("C9").Formula = "=C8+C3"
$ws.Range("D9").Formula = "=D8+D3"
$ws.Range("F9").Formula = "=F8+F3"
$ws.Range("G9").Formula = "=G8+G3"
$ws.Range("H9").Formula = "=H8+H3"

# --- Formatting ---
# Header row + phase/category labels in column B are bold.
$ws.Range("B2:H2").Font.Bold = $true
$ws.Range("B3:B9").Font.Bold = $true

# Merge the "Running costs" value down the phase rows and center it.
$ws.Range("E3:E7").Merge() | Out-Null
$ws.Range("E3:E7").HorizontalAlignment = -4108
$ws.Range("E3:E7").VerticalAlignment = -4108

# Stray formatted (underlined) empty cell left below the table.
$ws.Range("F16").Font.Underline = $true

# --- Column widths (best fit to content) ---
$ws.Columns.Item(2).ColumnWidth = 16.71
$ws.Columns.Item(3).ColumnWidth = 21.71
$ws.Columns.Item(4).ColumnWidth = 15
$ws.Columns.Item(5).ColumnWidth = 17.29
$ws.Columns.Item(6).ColumnWidth = 14.71
$ws.Columns.Item(7).ColumnWidth = 14.71
$ws.Columns.Item(8).ColumnWidth = 17.43

# --- Selection, matching the saved cursor position ---
$ws.Range("C12").Select() | Out-Null
